$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2 = "2025-10-17T07:09:42.340723"
    3 = "2025-10-17T07:09:42.340723"
    4 = "2025-10-17T07:09:42.341263"
    5 = "2025-10-17T07:09:42.341263"
    6 = "2025-10-17T07:09:42.341803"
    7 = "2025-10-17T07:09:42.341803"
    8 = "2025-10-17T07:09:42.341803"
    9 = "2025-10-17T07:09:42.342341"
    10 = "2025-10-17T07:09:42.342341"
    11 = "2025-10-17T07:09:42.342341"
    12 = "2025-10-17T07:09:42.342877"
    13 = "2025-10-17T07:09:42.342877"
    14 = "2025-10-17T07:09:42.343417"
    15 = "2025-10-17T07:09:42.343417"
    16 = "2025-10-17T07:09:42.343417"
    17 = "2025-10-17T07:09:42.343952"
    18 = "2025-10-17T07:09:42.343952"
    19 = "2025-10-17T07:09:42.344486"
    20 = "2025-10-17T07:09:42.344486"
    21 = "2025-10-17T07:09:42.344486"
    22 = "2025-10-17T07:09:42.345019"
    23 = "2025-10-17T07:09:42.345019"
    24 = "2025-10-17T07:09:42.345019"
    25 = "2025-10-17T07:09:42.345551"
    26 = "2025-10-17T07:09:42.345551"
    27 = "2025-10-17T07:09:42.346083"
    28 = "2025-10-17T07:09:42.346083"
    29 = "2025-10-17T07:09:42.346083"
    30 = "2025-10-17T07:09:42.346615"
    31 = "2025-10-17T07:09:42.346615"
    32 = "2025-10-17T07:09:42.347148"
    33 = "2025-10-17T07:09:42.347148"
    34 = "2025-10-17T07:09:42.347148"
    35 = "2025-10-17T07:09:42.347681"
    36 = "2025-10-17T07:09:42.347681"
    37 = "2025-10-17T07:09:42.347681"
    38 = "2025-10-17T07:09:42.348213"
    39 = "2025-10-17T07:09:42.348213"
    40 = "2025-10-17T07:09:42.348747"
    41 = "2025-10-17T07:09:42.348747"
    42 = "2025-10-17T07:09:42.348747"
    43 = "2025-10-17T07:09:42.349280"
    44 = "2025-10-17T07:09:42.349280"
    45 = "2025-10-17T07:09:42.349813"
    46 = "2025-10-17T07:09:42.411621"
    47 = "2025-10-17T07:09:42.411621"
    48 = "2025-10-17T07:09:42.412493"
    49 = "2025-10-17T07:09:42.412493"
    50 = "2025-10-17T07:09:42.413009"
    51 = "2025-10-17T07:09:42.413009"
    52 = "2025-10-17T07:09:42.413009"
    53 = "2025-10-17T07:09:42.413009"
    54 = "2025-10-17T07:09:42.414028"
    55 = "2025-10-17T07:09:42.414028"
    56 = "2025-10-17T07:09:42.414028"
    57 = "2025-10-17T07:09:42.414028"
    58 = "2025-10-17T07:09:42.415025"
    59 = "2025-10-17T07:09:42.415025"
    60 = "2025-10-17T07:09:42.415025"
    61 = "2025-10-17T07:09:42.415025"
    62 = "2025-10-17T07:09:42.415025"
    63 = "2025-10-17T07:09:42.416026"
    64 = "2025-10-17T07:09:42.416026"
    65 = "2025-10-17T07:09:42.416544"
    66 = "2025-10-17T07:09:42.416544"
    67 = "2025-10-17T07:09:42.416544"
    68 = "2025-10-17T07:09:42.417073"
    69 = "2025-10-17T07:09:42.417073"
    70 = "2025-10-17T07:09:42.417073"
    71 = "2025-10-17T07:09:42.417073"
    72 = "2025-10-17T07:09:42.417073"
    73 = "2025-10-17T07:09:42.418411"
    74 = "2025-10-17T07:09:42.418411"
    75 = "2025-10-17T07:09:42.498432"
    76 = "2025-10-17T07:09:42.498432"
    77 = "2025-10-17T07:09:42.498432"
    78 = "2025-10-17T07:09:42.498432"
    79 = "2025-10-17T07:09:42.499432"
    80 = "2025-10-17T07:09:42.499432"
    81 = "2025-10-17T07:09:42.499432"
    82 = "2025-10-17T07:09:42.499432"
    83 = "2025-10-17T07:09:42.499432"
    84 = "2025-10-17T07:09:42.500431"
    85 = "2025-10-17T07:09:42.500431"
    86 = "2025-10-17T07:09:42.500431"
    87 = "2025-10-17T07:09:42.500431"
    88 = "2025-10-17T07:09:42.500431"
    89 = "2025-10-17T07:09:42.500431"
    90 = "2025-10-17T07:09:42.500431"
    91 = "2025-10-17T07:09:42.501430"
    92 = "2025-10-17T07:09:42.501430"
    93 = "2025-10-17T07:09:42.501430"
    94 = "2025-10-17T07:09:42.501430"
    95 = "2025-10-17T07:09:42.501430"
    96 = "2025-10-17T07:09:42.502431"
    97 = "2025-10-17T07:09:42.502431"
    98 = "2025-10-17T07:09:42.502431"
    99 = "2025-10-17T07:09:42.502431"
    100 = "2025-10-17T07:09:42.502431"
    101 = "2025-10-17T07:09:42.502431"
    102 = "2025-10-17T07:09:42.503431"
    103 = "2025-10-17T07:09:42.556806"
    104 = "2025-10-17T07:09:42.556806"
    105 = "2025-10-17T07:09:42.556806"
    106 = "2025-10-17T07:09:42.556806"
    107 = "2025-10-17T07:09:42.556806"
    108 = "2025-10-17T07:09:42.556806"
    109 = "2025-10-17T07:09:42.556806"
    110 = "2025-10-17T07:09:42.556806"
    111 = "2025-10-17T07:09:42.556806"
    112 = "2025-10-17T07:09:42.556806"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item([int]$row, 26).Value = $timestamps[$row]
}

Write-Output "Updated timestamps for rows 2-112 in column Z"